$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain value updates (safe: will not be misinterpreted as numbers) ---
$ws.Range("D2").Value = '68.510.00'
$ws.Range("E2").Value = '  +2.16%  '
$ws.Range("D3").Value = '2.645.67'
$ws.Range("E3").Value = '  +1.65%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("E5").Value = '  +1.54%  '
$ws.Range("E6").Value = '  +3.07%  '
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("E8").Value = '  -0.32%  '
$ws.Range("D9").Value = '2.643.74'
$ws.Range("E9").Value = '  +1.63%  '
$ws.Range("E10").Value = '  +7.58%  '
$ws.Range("E11").Value = '  -0.39%  '
$ws.Range("E12").Value = '  +1.29%  '
$ws.Range("E13").Value = '  +1.19%  '
$ws.Range("E14").Value = '  +2.04%  '
$ws.Range("E15").Value = '  +4.21%  '
$ws.Range("D16").Value = '3.125.04'
$ws.Range("E16").Value = '  +1.77%  '
$ws.Range("D17").Value = '68.453.64'
$ws.Range("E17").Value = '  +2.24%  '
$ws.Range("D18").Value = '2.647.95'
$ws.Range("E18").Value = '  +1.82%  '
$ws.Range("E19").Value = '  +3.03%  '
$ws.Range("E20").Value = '  +1.04%  '
$ws.Range("E21").Value = '  +1.42%  '
$ws.Range("E22").Value = '  -0.15%  '
$ws.Range("E23").Value = '  +0.52%  '
$ws.Range("E24").Value = '  +2.77%  '
$ws.Range("E25").Value = '  +0.57%  '
$ws.Range("E26").Value = '  -0.18%  '
$ws.Range("E27").Value = '  +0.76%  '
$ws.Range("E28").Value = '  +6.56%  '
$ws.Range("D29").Value = '2.770.65'
$ws.Range("E30").Value = '  -0.07%  '
$ws.Range("E31").Value = '  -1.80%  '
$ws.Range("E32").Value = '  +4.20%  '
$ws.Range("E33").Value = '  +4.63%  '
$ws.Range("E34").Value = '  +2.50%  '
$ws.Range("E35").Value = '  +3.99%  '
$ws.Range("E36").Value = '  -0.01%  '
$ws.Range("E37").Value = '  +3.49%  '
$ws.Range("E38").Value = '  +2.66%  '
$ws.Range("E39").Value = '  +4.29%  '
$ws.Range("E40").Value = '  +1.35%  '
$ws.Range("B41").Value = 'PolygonEcosystemToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("E41").Value = '  +0.64%  '
$ws.Range("B42").Value = 'RenderToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("E42").Value = '  +3.54%  '
$ws.Range("E43").Value = '  +3.37%  '
$ws.Range("E44").Value = '  +3.81%  '
$ws.Range("E45").Value = '  +13.52%  '
$ws.Range("E46").Value = '  +0.04%  '
$ws.Range("E47").Value = '  -0.44%  '
$ws.Range("E48").Value = '  +2.75%  '
$ws.Range("E49").Value = '  +0.79%  '
$ws.Range("B50").Value = 'InjectiveProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("E50").Value = '  +2.42%  '
$ws.Range("B51").Value = 'Optimism'
$ws.Range("C51").Value = 'https://coinranking.com/coin/n1p-s_gm1+optimism-op'
$ws.Range("E51").Value = '  +1.77%  '

# --- Price cells that look numeric: force text storage via quote-prefix ---
# (mirrors typing e.g. '599.62 into Excel, which keeps it as literal text)
$ws.Range("D5").Value = "'599.62"
$ws.Range("D6").Value = "'154.62"
$ws.Range("D8").Value = "'0.545"
$ws.Range("D14").Value = "'27.84"
$ws.Range("D15").Value = "'0.0000187"
$ws.Range("D19").Value = "'11.38"
$ws.Range("D20").Value = "'367.39"
$ws.Range("D21").Value = "'7.44"
$ws.Range("D25").Value = "'73.28"
$ws.Range("D26").Value = "'0.998"
$ws.Range("D30").Value = "'0.999"
$ws.Range("D31").Value = "'574.62"
$ws.Range("D38").Value = "'159.53"
$ws.Range("D40").Value = "'19.22"
$ws.Range("D41").Value = "'0.367"
$ws.Range("D42").Value = "'5.39"
$ws.Range("D46").Value = "'0.999"
$ws.Range("D48").Value = "'157.28"
$ws.Range("D50").Value = "'21.96"
$ws.Range("D51").Value = "'1.71"
